$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (row 14) down to the two new rows
$ws.Range("A14:D14").Copy()
$ws.Range("A15:D16").PasteSpecial(-4122)  # xlPasteFormats

# Fill row 16 first so its description text is registered in the shared string
# table before row 15's text (matches original authoring order).
$ws.Range("A16").Value = 45586
$ws.Range("B16").Value = 1.5
$ws.Range("C16").Value = "Planung"
$ws.Range("D16").Value = "02_Product_Backlog.xlsx überarbeitet."

$ws.Range("A15").Value = 45586
$ws.Range("B15").Value = 0.5
$ws.Range("C15").Value = "Planung"
$ws.Range("D15").Value = "An Präsentation für Visions-Präsentation gearbeitet."

$ws.Range("D16").Select()
